# Add a 50th PSO run column to the results table.
#
# Previously, the last column (AZ) held the "Mean" summary statistic.
# A new run (Run 50) was performed, so:
#   - a new column is inserted at AZ for the "Run 50" data
#   - the old "Mean" column (formatting + data) shifts one column right to BA
#   - the "Mean" values are refreshed to reflect the new 50-run average

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a blank column at AZ; the existing "Mean" column (data + format)
# shifts right to BA automatically.
$ws.Columns("AZ:AZ").Insert()

# --- New "Run 50" column (AZ) ---
$ws.Range("AZ1").Value = "Run 50"

# Match the header formatting used by the other "Run n" header cells
# (bold font, centered/top alignment, thin border all around).
$ws.Range("AZ1").Font.Bold = $true
$ws.Range("AZ1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("AZ1").VerticalAlignment = -4160     # xlTop
$ws.Range("AZ1").Borders.LineStyle = 1

$runValue = 79.22083091
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = $runValue
}

# --- Refreshed "Mean" column (now BA) ---
$ws.Range("BA1").Value = "Mean"

$meanValue = 126.53409248
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 53).Value = $meanValue
}

Write-Output "Added Run 50 column and refreshed Mean column"
